$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.753.43"
$ws.Range("D3").Value = "2.446.26"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.62%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "2.445.36"
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("E11").Value = "  +2.56%  "
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("E13").Value = "  +1.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.23%  "
$ws.Range("E15").Value = "  +5.14%  "
$ws.Range("D16").Value = "2.889.69"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("D17").Value = "62.686.61"
$ws.Range("E17").Value = "  +3.31%  "
$ws.Range("D18").Value = "2.445.04"
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "330.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.86%  "
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("E23").Value = "  +7.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "646.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.53%  "
$ws.Range("E27").Value = "  +17.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.62%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0990"
$ws.Range("E29").Value = "  +5.15%  "
$ws.Range("B30").Value = "BabyDogeCoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D30").Value = "0.0₆0512"
$ws.Range("E30").Value = "  +84.30%  "
$ws.Range("E32").Value = "  +2.28%  "
$ws.Range("E33").Value = "  +6.72%  "
$ws.Range("E34").Value = "  +2.73%  "
$ws.Range("E35").Value = "  +4.26%  "
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("E38").Value = "  +3.05%  "
$ws.Range("E39").Value = "  +5.91%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.374"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "153.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.40%  "
$ws.Range("E43").Value = "  +7.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.17%  "
$ws.Range("E45").Value = "  +2.06%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "14.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +27.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "145.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.53%  "
$ws.Range("E49").Value = "  +3.48%  "
$ws.Range("E50").Value = "  +5.83%  "
$ws.Range("E51").Value = "  +2.27%  "
